$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = '[''FloorPlan_Train1_1'', ''FloorPlan_Train2_1'', ''FloorPlan_Train3_1'', ''FloorPlan_Train4_1'', ''FloorPlan_Train5_1'', ''FloorPlan_Train6_1'', ''FloorPlan_Train7_1'', ''FloorPlan_Train8_1'', ''FloorPlan_Train9_1'', ''FloorPlan_Train10_1'', ''FloorPlan_Train1_1'', ''FloorPlan_Train2_1'', ''FloorPlan_Train3_1'', ''FloorPlan_Train4_1'', ''FloorPlan_Train5_1'', ''FloorPlan_Train6_1'', ''FloorPlan_Train7_1'', ''FloorPlan_Train8_1'', ''FloorPlan_Train9_1'', ''FloorPlan_Train10_1'']'
$ws.Range("J2").Value = '[81, 93, 103, 114, 148, 194, 197, 197, 197, 197, 218, 220, 222, 240, 242, 293, 319, 319, 319, 319, 342]'
$ws.Range("K2").Value = '[0.086, 0.086, 0.098, 0.096, 0.094, 0.092, 0.092, 0.092, 0.092, 0.092, 0.092, 0.092, 0.092, 0.09, 0.09, 0.09, 0.1, 0.1, 0.1, 0.1, 0.1]'
$ws.Range("L2").Value = '[0, 30, 60, 90, 120, 150, 180, 210, 240, 270, 300, 330, 360, 390, 420, 450, 480, 510, 540, 570, 600]'
$ws.Range("M2").Value = '[0, 16.097865295410166, 38.60830733776093, 61.87577478885649, 83.7264050245285, 106.24110400676727, 113.78760421276093, 138.82562291622162, 151.63454074859618, 151.63454074859618, 160.57704050540923, 182.825630068779, 197.69045696258544, 219.51761140823365, 241.83314867019655, 267.2845475196839, 276.6690353274346, 303.59656013250355, 315.9621186614037, 315.9621186614037, 325.26231840848925]'
$ws.Range("N2").Value = '[0, 0.36, 0.86, 1.36, 1.93, 2.48, 3.06, 3.49, 3.98, 4.68, 5.11, 5.64, 6.21, 6.7, 7.61, 8.55, 9.45, 10.34, 11.35, 12.34, 13.03]'
$ws.Range("O2").Value = '[0, -10, -10, 0, -10, -10, -10, -10, 0, -10, 0, 0, 0, -10, -10, -20, 0, -20, -20, 0, -20, 0, -20, 0, -20, -20, 0, 0, 0, -20, 0, 0, -20, -20, -20, 0, -20, 0, 0, -20, -10, -20, -20, -20, -20, -20, -20, -10, -20, 0, -10, 0, 0, 0, 0, 0, 0, 0, -10, 0, 0, -10, 0, 0, -10, 0, 0, 0, -10, 0, 0, -10, 0, -10, 0, 0, 0, 0, 0, 0, -20]'

$ws.Range("I3").Value = '[''FloorPlan_Train1_1'', ''FloorPlan_Train2_1'', ''FloorPlan_Train3_1'', ''FloorPlan_Train4_1'', ''FloorPlan_Train5_1'', ''FloorPlan_Train6_1'', ''FloorPlan_Train7_1'', ''FloorPlan_Train8_1'', ''FloorPlan_Train9_1'', ''FloorPlan_Train10_1'', ''FloorPlan_Train1_1'', ''FloorPlan_Train2_1'', ''FloorPlan_Train3_1'', ''FloorPlan_Train4_1'', ''FloorPlan_Train5_1'', ''FloorPlan_Train6_1'', ''FloorPlan_Train7_1'', ''FloorPlan_Train8_1'', ''FloorPlan_Train9_1'', ''FloorPlan_Train10_1'']'
$ws.Range("J3").Value = '[81, 85, 100, 113, 143, 192, 213, 213, 213, 213, 237, 241, 244, 250, 275, 324, 347, 347, 347, 347, 393]'
$ws.Range("K3").Value = '[0.342, 0.352, 0.374, 0.37, 0.382, 0.382, 0.397, 0.397, 0.397, 0.397, 0.393, 0.398, 0.4, 0.4, 0.404, 0.404, 0.402, 0.402, 0.402, 0.402, 0.4]'
$ws.Range("L3").Value = '[0, 30, 60, 90, 120, 150, 180, 210, 240, 270, 300, 330, 360, 390, 420, 450, 480, 510, 540, 570, 600]'
$ws.Range("M3").Value = '[0, 14.544830417633067, 38.229873919487005, 60.96841380596162, 81.42432296276094, 105.40737187862398, 113.79237210750581, 134.84594135284428, 148.3509757995606, 148.3509757995606, 157.65117554664616, 178.0573615074158, 194.24799129962923, 219.1670227766037, 241.4825600385666, 266.9339588880539, 275.9897591352463, 302.9174396276474, 315.51484246253966, 315.51484246253966, 322.91660974025723]'
$ws.Range("N3").Value = '[0, 0.03, 0.08, 0.11, 0.16, 0.23, 0.28, 0.34, 0.39, 0.44, 0.5, 0.56, 0.63, 0.68, 0.75, 0.81, 0.9, 0.95, 1.01, 1.08, 1.14]'
$ws.Range("O3").Value = '[0, -10, -10, 0, -10, -10, -10, -10, 0, -10, 0, 0, 0, -10, -10, -20, 0, -20, -20, 0, 0, 0, -20, 0, -20, -20, 0, 0, 0, -20, 0, 0, -20, -20, -20, 0, -20, 0, 0, -20, -10, -20, -20, -20, 0, -20, -20, -20, -10, 0, -20, -10, 0, 0, 0, 0, 0, 0, 0, -10, -10, -10, 0, 0, -10, -10, 0, 0, 0, 0, 0, -10, 0, 0, -10, 0, -10, 0, 0, 0, -20]'

$ws.Range("I4").Value = '[''FloorPlan_Train1_1'', ''FloorPlan_Train2_1'', ''FloorPlan_Train3_1'', ''FloorPlan_Train4_1'', ''FloorPlan_Train5_1'', ''FloorPlan_Train6_1'', ''FloorPlan_Train7_1'', ''FloorPlan_Train8_1'', ''FloorPlan_Train9_1'', ''FloorPlan_Train10_1'', ''FloorPlan_Train11_1'', ''FloorPlan_Train11_2'', ''FloorPlan_Train11_3'', ''FloorPlan_Train11_4'', ''FloorPlan_Train11_5'', ''FloorPlan_Train12_1'', ''FloorPlan_Train12_2'', ''FloorPlan_Train12_3'', ''FloorPlan_Train12_4'', ''FloorPlan_Train12_5'']'
$ws.Range("J4").Value = '[81, 84, 94, 102, 142, 189, 190, 190, 190, 190, 210, 210, 210, 210, 210, 215, 218, 259, 306, 306, 329]'
$ws.Range("K4").Value = '[0.086, 0.086, 0.092, 0.092, 0.094, 0.103, 0.103, 0.103, 0.103, 0.103, 0.103, 0.103, 0.103, 0.103, 0.103, 0.103, 0.103, 0.109, 0.107, 0.107, 0.105]'
$ws.Range("L4").Value = '[0, 30, 60, 90, 120, 150, 180, 210, 240, 270, 300, 330, 360, 390, 420, 450, 480, 510, 540, 570, 600]'
$ws.Range("M4").Value = '[0, 15.530802631378185, 37.09476878643035, 57.24029161930084, 77.69620077610017, 100.7003497362137, 107.74374992847444, 132.30935318470003, 145.43205032348635, 145.43205032348635, 154.19570008516314, 175.495841896534, 192.57426539659502, 211.8608514428139, 235.58418039083483, 258.6813417553902, 275.4413020730019, 292.6162026047707, 307.3982270836831, 307.8179270863534, 328.02542659044275]'
$ws.Range("N4").Value = '[0, 0.31, 0.67, 1.23, 2.04, 2.73, 3.32, 3.81, 4.39, 4.78, 5.2, 5.58, 6.0, 6.52, 7.09, 7.61, 8.33, 9.29, 10.2, 11.25, 12.3]'
$ws.Range("O4").Value = '[-10, 0, 0, -10, -20, 0, 0, -10, -20, 0, -10, 0, -20, -20, 0, 0, 0, 0, 0, -10, -10, 0, 0, 0, 0, 0, -20, -10, -10, 0, -10, -20, 0, 0, 0, -10, -20, -20, 0, 0, 0, 0, 0, 0, -20, 0, -20, -10, 0, -20, 0, 0, -10, 0, 0, 0, 0, -20, -20, 0, -20, -20, -10, 0, -20, 0, 0, -10, -10, -20, -10, -10, 0, 0, -20, -20, -10, -10, -20, -10, 0]'

$ws.Range("I5").Value = '[''FloorPlan_Train1_1'', ''FloorPlan_Train2_1'', ''FloorPlan_Train3_1'', ''FloorPlan_Train4_1'', ''FloorPlan_Train5_1'', ''FloorPlan_Train6_1'', ''FloorPlan_Train7_1'', ''FloorPlan_Train8_1'', ''FloorPlan_Train9_1'', ''FloorPlan_Train10_1'', ''FloorPlan_Train1_4'', ''FloorPlan_Train2_4'', ''FloorPlan_Train3_4'', ''FloorPlan_Train4_4'', ''FloorPlan_Train5_4'', ''FloorPlan_Train6_4'', ''FloorPlan_Train7_4'', ''FloorPlan_Train8_4'', ''FloorPlan_Train9_4'', ''FloorPlan_Train10_4'']'
$ws.Range("J5").Value = '[81, 87, 98, 104, 144, 192, 214, 214, 214, 214, 223, 223, 227, 250, 252, 252, 252, 252, 252, 252, 252]'
$ws.Range("K5").Value = '[0.342, 0.365, 0.389, 0.385, 0.383, 0.383, 0.393, 0.393, 0.393, 0.393, 0.397, 0.397, 0.404, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41, 0.41]'
$ws.Range("L5").Value = '[0, 30, 60, 90, 120, 150, 180, 210, 240, 270, 300, 330, 360, 390, 420, 450, 480, 510, 540, 570, 600]'
$ws.Range("M5").Value = '[0, 16.67876248359681, 38.05118577480316, 60.16182434558867, 80.61773350238799, 104.11133244037627, 112.49633266925811, 138.0067668199539, 151.75726034641264, 151.75726034641264, 173.77884027361867, 201.77658314108845, 213.1931896626949, 237.1615375936031, 250.92134079337117, 273.1265073239803, 299.4534540593624, 320.50671071410176, 340.2187743127346, 361.99377297759054, 368.4016011893749]'
$ws.Range("N5").Value = '[0, 0.04, 0.07, 0.1, 0.13, 0.18, 0.23, 0.29, 0.35, 0.42, 0.46, 0.51, 0.6, 0.69, 0.75, 0.81, 0.86, 0.91, 0.97, 1.06, 1.13]'
$ws.Range("O5").Value = '[-10, 0, -10, -20, -10, 0, -20, 0, 0, -10, -20, 0, 0, 0, -10, 0, 0, 0, 0, 0, -20, 0, -20, 0, -20, 0, 0, -20, -20, 0, 0, -20, 0, 0, 0, -10, -10, -20, 0, 0, 0, -10, -10, -20, 0, -20, 0, -20, 0, 0, -10, -10, 0, 0, -20, 0, -10, 0, 0, 0, 0, -10, -20, 0, -20, -10, -10, -10, -10, -20, 0, 0, 0, -10, -10, 0, -20, -20, -20, 0, -10]'
